$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.571.06"
$ws.Range("E2").Value = "  +3.78%  "

# Row 3
$ws.Range("D3").Value = "2.647.62"
$ws.Range("E3").Value = "  +2.55%  "

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "'568.80"
$ws.Range("E5").Value = "  +5.92%  "

# Row 6
$ws.Range("D6").Value = "'146.04"
$ws.Range("E6").Value = "  +2.08%  "

# Row 7
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  +4.25%  "

# Row 9
$ws.Range("D9").Value = "2.647.29"
$ws.Range("E9").Value = "  +2.43%  "

# Row 10
$ws.Range("D10").Value = "'6.83"
$ws.Range("E10").Value = "  +1.31%  "

# Row 11
$ws.Range("E11").Value = "  +5.07%  "

# Row 12
$ws.Range("D12").Value = "'0.151"
$ws.Range("E12").Value = "  +9.56%  "

# Row 13
$ws.Range("E13").Value = "  +4.28%  "

# Row 14
$ws.Range("D14").Value = "3.107.56"
$ws.Range("E14").Value = "  +2.41%  "

# Row 15
$ws.Range("D15").Value = "60.506.38"
$ws.Range("E15").Value = "  +3.83%  "

# Row 16
$ws.Range("D16").Value = "'22.06"
$ws.Range("E16").Value = "  +6.73%  "

# Row 17
$ws.Range("E17").Value = "  +3.04%  "

# Row 18
$ws.Range("D18").Value = "2.642.67"
$ws.Range("E18").Value = "  +1.82%  "

# Row 19
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  +2.20%  "

# Row 20
$ws.Range("D20").Value = "'340.96"
$ws.Range("E20").Value = "  +1.73%  "

# Row 21
$ws.Range("E21").Value = "  +3.95%  "

# Row 22
$ws.Range("E22").Value = "  +3.67%  "

# Row 23
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
$ws.Range("D24").Value = "'66.11"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25
$ws.Range("D25").Value = "'0.447"
$ws.Range("E25").Value = "  +6.87%  "

# Row 26
$ws.Range("D26").Value = "'0.164"
$ws.Range("E26").Value = "  +4.87%  "

# Row 27
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.28%  "

# Row 28
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  +4.68%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0800"
$ws.Range("E29").Value = "  +9.96%  "

# Row 30
$ws.Range("E30").Value = "  -0.08%  "

# Row 31
$ws.Range("E31").Value = "  +4.51%  "

# Row 32
$ws.Range("D32").Value = "'6.13"
$ws.Range("E32").Value = "  +4.08%  "

# Row 33
$ws.Range("D33").Value = "'158.36"
$ws.Range("E33").Value = "  +2.75%  "

# Row 34
$ws.Range("E34").Value = "  +1.82%  "

# Row 35
$ws.Range("E35").Value = "  +5.71%  "

# Row 36
$ws.Range("E36").Value = "  +5.62%  "

# Row 37
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "'0.882"
$ws.Range("E37").Value = "  +7.02%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.885"
$ws.Range("E38").Value = "  +8.22%  "

# Row 39
$ws.Range("D39").Value = "'37.59"
$ws.Range("E39").Value = "  +1.75%  "

# Row 40
$ws.Range("D40").Value = "'1.51"
$ws.Range("E40").Value = "  +7.15%  "

# Row 41
$ws.Range("D41").Value = "'300.47"
$ws.Range("E41").Value = "  +6.12%  "

# Row 42
$ws.Range("E42").Value = "  +1.93%  "

# Row 43
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("E44").Value = "  +4.82%  "

# Row 45
$ws.Range("D45").Value = "'0.602"
$ws.Range("E45").Value = "  +2.23%  "

# Row 46
$ws.Range("E46").Value = "  +2.20%  "

# Row 47
$ws.Range("D47").Value = "'19.41"
$ws.Range("E47").Value = "  +5.46%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'127.84"
$ws.Range("E48").Value = "  +16.36%  "

# Row 49
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.69"
$ws.Range("E49").Value = "  +0.52%  "

# Row 50
$ws.Range("E50").Value = "  +4.06%  "

# Row 51
$ws.Range("D51").Value = "'4.66"
$ws.Range("E51").Value = "  +6.47%  "
